$wb = $excel.ActiveWorkbook

$summary = $wb.Worksheets.Item("总计")
$q3 = $wb.Worksheets.Item("2022-Q3")

# 1) Duplicate the existing "2022-Q3" sheet right after itself so the old
#    Q3 figures (and their original formatting) survive on their own tab.
$q3.Copy($null, $q3)
$q3dup = $wb.Worksheets.Item(3)

# 2) Free up the "2022-Q3" name on the original sheet, then let the
#    duplicate take it back - the duplicate now permanently holds the
#    untouched Q3 data/style, the original becomes the Q4 sheet.
$q3.Name = "2022-Q4-tmp"
$q3dup.Name = "2022-Q3"
$q3.Name = "2022-Q4"

# 3) Refresh the figures on the (now) "2022-Q4" sheet with this quarter's
#    numbers.
$q3.Range("D2").Value = 4.63
$q3.Range("E2").Value = 24.21
$q3.Range("F2").Value = 4.35
$q3.Range("G2").Value = 0.2014

# 4) Re-style the "2022-Q4" header row + A2 to match the "总计" sheet's
#    header formatting.
$summary.Range("B1:D1").Copy()
$q3.Range("B1").PasteSpecial(-4122)
$summary.Range("B1").Copy()
$q3.Range("E1:H1").PasteSpecial(-4122)
$summary.Range("A2").Copy()
$q3.Range("A2").PasteSpecial(-4122)

# 5) Match the "总计" sheet's page margins on "2022-Q4".
$q3.PageSetup.LeftMargin = 54
$q3.PageSetup.RightMargin = 54
$q3.PageSetup.TopMargin = 72
$q3.PageSetup.BottomMargin = 72
$q3.PageSetup.HeaderMargin = 36
$q3.PageSetup.FooterMargin = 36

# 6) Update the "总计" roll-up sheet: point row 2 at 2022-Q4 and append a
#    new row 3 carrying the prior 2022-Q3 totals.
$summary.Range("B2").Value = "2022-Q4"
$summary.Range("D2").Value = 0.2
$summary.Range("A3").Value = 1
$summary.Range("B3").Value = "2022-Q3"
$summary.Range("C3").Value = 1
$summary.Range("D3").Value = 0.21
$summary.Range("A2").Copy()
$summary.Range("A3").PasteSpecial(-4122)
